$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("E6:E13")
$rng.Select()
$rng.ClearContents()
